$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12, which shifts the existing rows 12-115
# (and all their data) down to rows 13-116, growing the used range from
# A1:R115 to A1:R116.
$ws.Rows("12:12").Insert()

# Populate the newly inserted row 12 with its data (a new weekly price
# record for "Poroto verde").
$ws.Cells.Item(12, 1).Value  = 7
$ws.Cells.Item(12, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(12, 3).Value  = "Ñuble"
$ws.Cells.Item(12, 4).Value  = 44950
$ws.Cells.Item(12, 5).Value  = 16
$ws.Cells.Item(12, 6).Value  = 100112031
$ws.Cells.Item(12, 7).Value  = "Poroto verde"
$ws.Cells.Item(12, 8).Value  = "Sin especificar"
$ws.Cells.Item(12, 9).Value  = "Primera"
$ws.Cells.Item(12, 10).Value = 50
$ws.Cells.Item(12, 11).Value = 28000
$ws.Cells.Item(12, 12).Value = 28000
$ws.Cells.Item(12, 13).Value = 28000
$ws.Cells.Item(12, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(12, 15).Value = "Región del Maule"
$ws.Cells.Item(12, 16).Value = 1120
$ws.Cells.Item(12, 17).Value = 25
$ws.Cells.Item(12, 18).Value = "Hortaliza"
